$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A -- this shifts the existing
# Name/Surname columns from A/B to B/C and carries their widths along.
$ws.Columns.Item(1).Insert()

# Populate the new "Id" column: header + row numbers 1, 2, 3
$ws.Range("A1").Value = "Id"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

# Give the new Id column a narrow width (closest reachable approximation
# of the real-Excel auto-fit width for this content)
$ws.Columns.Item(1).ColumnWidth = 1.35

# Leave the sheet with the whole-grid selected, as happens after an
# insert-column operation performed from the column header context menu
$ws.Cells.Select() | Out-Null
